$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 95 (old rows 95:98 shift down to 99:102)
$ws.Rows("95:98").Insert()

# Row 95: EtudeSkiniHarmonie1
$ws.Range("C95").Value = "EtudeSkiniHarmonie1.csv"
$ws.Range("F95").Value = "EtudeSkiniHarmonie1"
$ws.Range("G95").Value = "EtudeSkiniHarmonie1.xml"

# Row 96: EtudeSkiniHarmonie2et3
$ws.Range("C96").Value = "EtudeSkiniHarmonie2.csv"
$ws.Range("F96").Value = "EtudeSkiniHarmonie2et3"
$ws.Range("G96").Value = "EtudeSkiniHarmonie2et3.xml"

# Row 97: EtudeSkiniHarmonie4
$ws.Range("C97").Value = "EtudeSkiniHarmonie4.csv"
$ws.Range("F97").Value = "EtudeSkiniHarmonie4"
$ws.Range("G97").Value = "EtudeSkiniHarmonie4.xml"

# Row 98: EtudeSkiniHarmonie5
$ws.Range("C98").Value = "EtudeSkiniHarmonie5.csv"
$ws.Range("F98").Value = "EtudeSkiniHarmonie5"
$ws.Range("G98").Value = "EtudeSkiniHarmonie5.xml"

# Update selection to match final cursor position
$ws.Range("C94").Select()
